$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The WBS table shrank from 55 data+header rows (A6:D60) to 44 (A6:D49).
# Remove the now-unused trailing rows first (this also auto-shrinks the
# Tabelle1 ListObject range from A6:D60 to A6:D49).
$ws.Range("A50:D60").EntireRow.Delete()

# Rewrite the table header (row 6) and all data rows (7-49) with the updated
# WBS content: "Main Character" renamed to "Hammer Man", left/right movement
# rows merged into single "movement" rows, "spawn waves" renamed to
# "scriptable spawn waves", and the generic Spalte1..4 headers replaced with
# their real names.
$ws.Cells.Item(6, 1).Value = 'Category'
$ws.Cells.Item(6, 2).Value = 'Goal'
$ws.Cells.Item(6, 3).Value = 'Task'
$ws.Cells.Item(6, 4).Value = 'Department'

$ws.Cells.Item(7, 1).Value = 'Hammer Man'
$ws.Cells.Item(7, 2).Value = 'movement'
$ws.Cells.Item(7, 3).Value = 'implementation'
$ws.Cells.Item(7, 4).Value = 'Code'

$ws.Cells.Item(8, 1).Value = 'Hammer Man'
$ws.Cells.Item(8, 2).Value = 'movement'
$ws.Cells.Item(8, 3).Value = 'sprite of Hammer Man moving with hammer over head'
$ws.Cells.Item(8, 4).Value = 'Art'

$ws.Cells.Item(9, 1).Value = 'Hammer Man'
$ws.Cells.Item(9, 2).Value = 'movement'
$ws.Cells.Item(9, 3).Value = 'sprite of Hammer Man moving with hammer at high mid height'
$ws.Cells.Item(9, 4).Value = 'Art'

$ws.Cells.Item(10, 1).Value = 'Hammer Man'
$ws.Cells.Item(10, 2).Value = 'movement'
$ws.Cells.Item(10, 3).Value = 'sprite of Hammer Man moving with hammer at low mid height'
$ws.Cells.Item(10, 4).Value = 'Art'

$ws.Cells.Item(11, 1).Value = 'Hammer Man'
$ws.Cells.Item(11, 2).Value = 'movement'
$ws.Cells.Item(11, 3).Value = 'sprite of Hammer Man moving with hammer at low height'
$ws.Cells.Item(11, 4).Value = 'Art'

$ws.Cells.Item(12, 1).Value = 'Hammer Man'
$ws.Cells.Item(12, 2).Value = 'movement'
$ws.Cells.Item(12, 3).Value = 'animation of Hammer Man moving'
$ws.Cells.Item(12, 4).Value = 'Art'

$ws.Cells.Item(13, 1).Value = 'Hammer Man'
$ws.Cells.Item(13, 2).Value = 'movement climb'
$ws.Cells.Item(13, 3).Value = 'implementation'
$ws.Cells.Item(13, 4).Value = 'Code'

$ws.Cells.Item(14, 1).Value = 'Hammer Man'
$ws.Cells.Item(14, 2).Value = 'movement climb'
$ws.Cells.Item(14, 3).Value = 'sprite of Hammer Man climbing up (right foot up)'
$ws.Cells.Item(14, 4).Value = 'Art'

$ws.Cells.Item(15, 1).Value = 'Hammer Man'
$ws.Cells.Item(15, 2).Value = 'movement climb'
$ws.Cells.Item(15, 3).Value = 'sprite of Hammer Man climbing up (right foot down)'
$ws.Cells.Item(15, 4).Value = 'Art'

$ws.Cells.Item(16, 1).Value = 'Hammer Man'
$ws.Cells.Item(16, 2).Value = 'hammer attack'
$ws.Cells.Item(16, 3).Value = 'implementation (hit box moves)'
$ws.Cells.Item(16, 4).Value = 'Code'

$ws.Cells.Item(17, 1).Value = 'Hammer Man'
$ws.Cells.Item(17, 2).Value = 'hammer attack'
$ws.Cells.Item(17, 3).Value = 'hammer hit enemy sound (umpf!)'
$ws.Cells.Item(17, 4).Value = 'Audio'

$ws.Cells.Item(18, 1).Value = 'Enemy'
$ws.Cells.Item(18, 2).Value = 'movement'
$ws.Cells.Item(18, 3).Value = 'implementation'
$ws.Cells.Item(18, 4).Value = 'Code'

$ws.Cells.Item(19, 1).Value = 'Enemy'
$ws.Cells.Item(19, 2).Value = 'movement'
$ws.Cells.Item(19, 3).Value = 'enemy moves (left foot forward)'
$ws.Cells.Item(19, 4).Value = 'Art'

$ws.Cells.Item(20, 1).Value = 'Enemy'
$ws.Cells.Item(20, 2).Value = 'movement'
$ws.Cells.Item(20, 3).Value = 'enemy moves (right foot forward) '
$ws.Cells.Item(20, 4).Value = 'Art'

$ws.Cells.Item(21, 1).Value = 'Enemy'
$ws.Cells.Item(21, 2).Value = 'movement'
$ws.Cells.Item(21, 3).Value = 'animation of enemy moving'
$ws.Cells.Item(21, 4).Value = 'Art'

$ws.Cells.Item(22, 1).Value = 'Enemy'
$ws.Cells.Item(22, 2).Value = 'turn at wall'
$ws.Cells.Item(22, 3).Value = 'implementation'
$ws.Cells.Item(22, 4).Value = 'Code'

$ws.Cells.Item(23, 1).Value = 'Enemy'
$ws.Cells.Item(23, 2).Value = 'randomly change direction after fall'
$ws.Cells.Item(23, 3).Value = 'implementation'
$ws.Cells.Item(23, 4).Value = 'Code'

$ws.Cells.Item(24, 1).Value = 'Enemy'
$ws.Cells.Item(24, 2).Value = 'spawn and randomly walk left or right'
$ws.Cells.Item(24, 3).Value = 'implementation'
$ws.Cells.Item(24, 4).Value = 'Code'

$ws.Cells.Item(25, 1).Value = 'Enemy'
$ws.Cells.Item(25, 2).Value = 'exit level'
$ws.Cells.Item(25, 3).Value = 'implementation'
$ws.Cells.Item(25, 4).Value = 'Code'

$ws.Cells.Item(26, 1).Value = 'Enemy'
$ws.Cells.Item(26, 2).Value = 'scriptable spawn waves'
$ws.Cells.Item(26, 3).Value = 'implementation'
$ws.Cells.Item(26, 4).Value = 'Code'

$ws.Cells.Item(27, 1).Value = 'Enemy'
$ws.Cells.Item(27, 2).Value = 'counter for how many of the wave are left'
$ws.Cells.Item(27, 3).Value = 'implementation'
$ws.Cells.Item(27, 4).Value = 'Code'

$ws.Cells.Item(28, 1).Value = 'Enemy'
$ws.Cells.Item(28, 2).Value = 'counter for how many of the wave are left'
$ws.Cells.Item(28, 3).Value = 'bar representing left over enemies'
$ws.Cells.Item(28, 4).Value = 'Art'

$ws.Cells.Item(29, 1).Value = 'Enemy'
$ws.Cells.Item(29, 2).Value = 'enemy deals damage on exit'
$ws.Cells.Item(29, 3).Value = 'implementation'
$ws.Cells.Item(29, 4).Value = 'Code'

$ws.Cells.Item(30, 1).Value = 'Enemy'
$ws.Cells.Item(30, 2).Value = 'die on hit'
$ws.Cells.Item(30, 3).Value = 'implementation'
$ws.Cells.Item(30, 4).Value = 'Code'

$ws.Cells.Item(31, 1).Value = 'Level'
$ws.Cells.Item(31, 2).Value = 'Layout'
$ws.Cells.Item(31, 3).Value = 'design level layout'
$ws.Cells.Item(31, 4).Value = 'GD'

$ws.Cells.Item(32, 1).Value = 'Level'
$ws.Cells.Item(32, 2).Value = 'Layout'
$ws.Cells.Item(32, 3).Value = 'apply textures in level'
$ws.Cells.Item(32, 4).Value = 'GD'

$ws.Cells.Item(33, 1).Value = 'Level'
$ws.Cells.Item(33, 2).Value = 'platforms'
$ws.Cells.Item(33, 3).Value = 'sprite of a platform'
$ws.Cells.Item(33, 4).Value = 'Art'

$ws.Cells.Item(34, 1).Value = 'Level'
$ws.Cells.Item(34, 2).Value = 'ladders'
$ws.Cells.Item(34, 3).Value = 'sprite of a ladder'
$ws.Cells.Item(34, 4).Value = 'Art'

$ws.Cells.Item(35, 1).Value = 'Level'
$ws.Cells.Item(35, 2).Value = 'spawn point'
$ws.Cells.Item(35, 3).Value = 'sprite of a spawn point'
$ws.Cells.Item(35, 4).Value = 'Art'

$ws.Cells.Item(36, 1).Value = 'Level'
$ws.Cells.Item(36, 2).Value = 'exit point'
$ws.Cells.Item(36, 3).Value = 'sprite of a exit point'
$ws.Cells.Item(36, 4).Value = 'Art'

$ws.Cells.Item(37, 1).Value = 'Menu Screen'
$ws.Cells.Item(37, 2).Value = 'start game button'
$ws.Cells.Item(37, 3).Value = 'implementation'
$ws.Cells.Item(37, 4).Value = 'Code'

$ws.Cells.Item(38, 1).Value = 'Menu Screen'
$ws.Cells.Item(38, 2).Value = 'start game button'
$ws.Cells.Item(38, 3).Value = 'button design'
$ws.Cells.Item(38, 4).Value = 'Art'

$ws.Cells.Item(39, 1).Value = 'Menu Screen'
$ws.Cells.Item(39, 2).Value = 'exit game button'
$ws.Cells.Item(39, 3).Value = 'implementation'
$ws.Cells.Item(39, 4).Value = 'Code'

$ws.Cells.Item(40, 1).Value = 'Menu Screen'
$ws.Cells.Item(40, 2).Value = 'exit game button'
$ws.Cells.Item(40, 3).Value = 'button design'
$ws.Cells.Item(40, 4).Value = 'Art'

$ws.Cells.Item(41, 1).Value = 'Menu Screen'
$ws.Cells.Item(41, 2).Value = 'credits button'
$ws.Cells.Item(41, 3).Value = 'implementation'
$ws.Cells.Item(41, 4).Value = 'Code'

$ws.Cells.Item(42, 1).Value = 'Menu Screen'
$ws.Cells.Item(42, 2).Value = 'credits button'
$ws.Cells.Item(42, 3).Value = 'button design'
$ws.Cells.Item(42, 4).Value = 'Art'

$ws.Cells.Item(43, 1).Value = 'Credits Screen'
$ws.Cells.Item(43, 2).Value = 'credits'
$ws.Cells.Item(43, 3).Value = 'credit design'
$ws.Cells.Item(43, 4).Value = 'Art'

$ws.Cells.Item(44, 1).Value = 'Game Over Screen'
$ws.Cells.Item(44, 2).Value = 'retry game button'
$ws.Cells.Item(44, 3).Value = 'implementation'
$ws.Cells.Item(44, 4).Value = 'Code'

$ws.Cells.Item(45, 1).Value = 'Game Over Screen'
$ws.Cells.Item(45, 2).Value = 'retry game button'
$ws.Cells.Item(45, 3).Value = 'button design'
$ws.Cells.Item(45, 4).Value = 'Art'

$ws.Cells.Item(46, 1).Value = 'Game Over Screen'
$ws.Cells.Item(46, 2).Value = 'exit to menu screen'
$ws.Cells.Item(46, 3).Value = 'implementation'
$ws.Cells.Item(46, 4).Value = 'Code'

$ws.Cells.Item(47, 1).Value = 'Game Over Screen'
$ws.Cells.Item(47, 2).Value = 'exit to menu screen'
$ws.Cells.Item(47, 3).Value = 'button design'
$ws.Cells.Item(47, 4).Value = 'Art'

$ws.Cells.Item(48, 1).Value = 'UI'
$ws.Cells.Item(48, 2).Value = 'health bar'
$ws.Cells.Item(48, 3).Value = 'implementation'
$ws.Cells.Item(48, 4).Value = 'Code'

$ws.Cells.Item(49, 1).Value = 'UI'
$ws.Cells.Item(49, 2).Value = 'health bar'
$ws.Cells.Item(49, 3).Value = 'health bar design'
$ws.Cells.Item(49, 4).Value = 'Art'

# Update the view to match the author's final scroll/selection state.
$ws.Range("B26").Select()
